# Update test for correctly testing time value as string.
#
# Cells F7 and G7 previously held the numeric TimeSpan value
# 1.10538194444444 (formatted with the "[h]:mm:ss" style, style index 5)
# representing 1 day, 2 hours, 31 minutes, 45 seconds.
#
# They should instead hold the literal text "1.02:31:45" (the .NET
# TimeSpan.ToString() representation of that same duration), stored as a
# shared string, and revert to the default (unstyled / style index 0)
# cell format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timeSpanText = "1.02:31:45"

$cellF7 = $ws.Range("F7")
$cellF7.NumberFormat = ""
$cellF7.Value = $timeSpanText

$cellG7 = $ws.Range("G7")
$cellG7.NumberFormat = ""
$cellG7.Value = $timeSpanText
